$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Gamma1F"

# Tiny last-bit precision corrections from recomputation (row 13)
$ws.Cells.Item(13, 3).Value = 0.9902816289678632
$ws.Cells.Item(13, 4).Value = 0.9916090644196027
$ws.Cells.Item(13, 6).Value = 0.9902816289678632
$ws.Cells.Item(13, 10).Value = 0.9916090644196027
$ws.Cells.Item(13, 11).Value = 0.9913549366812358

# Tiny last-bit precision corrections from recomputation (row 15)
$ws.Cells.Item(15, 3).Value = 0.9888161938134855
$ws.Cells.Item(15, 6).Value = 0.9888161938134855

# New row 16 data
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 0.9997970327128718
$ws.Cells.Item(16, 4).Value = 0.9697660149721196
$ws.Cells.Item(16, 5).Value = 0.9997385364559017
$ws.Cells.Item(16, 6).Value = 0.9997970327128718
$ws.Cells.Item(16, 7).Value = 0.9699415614231334
$ws.Cells.Item(16, 8).Value = 1.001031999558077
$ws.Cells.Item(16, 9).Value = 0.9941468951873086
$ws.Cells.Item(16, 10).Value = 0.9697660149721196
$ws.Cells.Item(16, 11).Value = 0.9847522757140106
$ws.Cells.Item(16, 12).Value = 0.9922746542134413
$ws.Cells.Item(16, 13).Value = 0.9890703400515687
